# Auto-generated edit script: apply scheduled-runner market-data refresh
# to the per-sheet Leve profit tables (columns H:N) as captured by the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Cells.Item(6, 8).Value = 200.25  # H6: 228.71428 -> 200.25
$ws.Cells.Item(6, 9).Value = 228.71428  # I6: 266.66666 -> 228.71428
$ws.Cells.Item(6, 11).Value = 686.14284  # K6: 799.9999799999999 -> 686.14284
$ws.Cells.Item(6, 13).Value = -574.14284  # M6: -687.9999799999999 -> -574.14284

# Row 19
$ws.Cells.Item(19, 8).Value = 2500  # H19: 1990.8 -> 2500
$ws.Cells.Item(19, 9).Value = 0  # I19: 1499.5 -> 0
$ws.Cells.Item(19, 10).Value = 2500  # J19: 2318.3333 -> 2500
$ws.Cells.Item(19, 11).Value = 0  # K19: 1499.5 -> 0
$ws.Cells.Item(19, 12).Value = 2500  # L19: 2318.3333 -> 2500
$ws.Cells.Item(19, 13).ClearContents()  # M19: -1324.5 -> (removed)
$ws.Cells.Item(19, 14).Value = -2850  # N19: -2668.3333 -> -2850

# Row 33
$ws.Cells.Item(33, 8).Value = 3474248.2  # H33: 3474248.5 -> 3474248.2

# Row 43
$ws.Cells.Item(43, 8).Value = 3776.625  # H43: 3894.6 -> 3776.625
$ws.Cells.Item(43, 9).Value = 3615  # I43: 3998 -> 3615
$ws.Cells.Item(43, 10).Value = 3873.6  # J43: 3868.75 -> 3873.6
$ws.Cells.Item(43, 11).Value = 3615  # K43: 3998 -> 3615
$ws.Cells.Item(43, 12).Value = 3873.6  # L43: 3868.75 -> 3873.6
$ws.Cells.Item(43, 13).Value = -3546  # M43: -3929 -> -3546
$ws.Cells.Item(43, 14).Value = -4011.6  # N43: -4006.75 -> -4011.6

# Row 57
$ws.Cells.Item(57, 8).Value = 68850  # H57: 74140.336 -> 68850
$ws.Cells.Item(57, 10).Value = 68850  # J57: 74140.336 -> 68850
$ws.Cells.Item(57, 12).Value = 206550  # L57: 222421.008 -> 206550
$ws.Cells.Item(57, 14).Value = -207548  # N57: -223419.008 -> -207548

# Row 111
$ws.Cells.Item(111, 8).Value = 4118.75  # H111: 4137.4287 -> 4118.75
$ws.Cells.Item(111, 9).Value = 3994.25  # I111: 3994.5 -> 3994.25
$ws.Cells.Item(111, 10).Value = 4243.25  # J111: 4328 -> 4243.25
$ws.Cells.Item(111, 11).Value = 11982.75  # K111: 11983.5 -> 11982.75
$ws.Cells.Item(111, 12).Value = 12729.75  # L111: 12984 -> 12729.75
$ws.Cells.Item(111, 13).Value = -8915.75  # M111: -8916.5 -> -8915.75
$ws.Cells.Item(111, 14).Value = -18863.75  # N111: -19118 -> -18863.75

# Row 138
$ws.Cells.Item(138, 8).Value = 1182.2  # H138: 1000.38464 -> 1182.2
$ws.Cells.Item(138, 9).Value = 1182.2  # I138: 1000.38464 -> 1182.2
$ws.Cells.Item(138, 11).Value = 3546.6  # K138: 3001.15392 -> 3546.6
$ws.Cells.Item(138, 13).Value = 1593.4  # M138: 2138.84608 -> 1593.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 1384.4894  # H32: 1384.5957 -> 1384.4894
$ws.Cells.Item(32, 9).Value = 1384.4894  # I32: 1384.5957 -> 1384.4894
$ws.Cells.Item(32, 11).Value = 1384.4894  # K32: 1384.5957 -> 1384.4894
$ws.Cells.Item(32, 13).Value = -1097.4894  # M32: -1097.5957 -> -1097.4894

# Row 74
$ws.Cells.Item(74, 8).Value = 6615041.5  # H74: 4210504 -> 6615041.5
$ws.Cells.Item(74, 9).Value = 3704903.5  # I74: 2059418.8 -> 3704903.5
$ws.Cells.Item(74, 10).Value = 13890387  # J74: 13890388 -> 13890387
$ws.Cells.Item(74, 11).Value = 3704903.5  # K74: 2059418.8 -> 3704903.5
$ws.Cells.Item(74, 12).Value = 13890387  # L74: 13890388 -> 13890387
$ws.Cells.Item(74, 13).Value = -3704029.5  # M74: -2058544.8 -> -3704029.5
$ws.Cells.Item(74, 14).Value = -13892135  # N74: -13892136 -> -13892135

# Row 77
$ws.Cells.Item(77, 8).Value = 6615041.5  # H77: 4210504 -> 6615041.5
$ws.Cells.Item(77, 9).Value = 3704903.5  # I77: 2059418.8 -> 3704903.5
$ws.Cells.Item(77, 10).Value = 13890387  # J77: 13890388 -> 13890387
$ws.Cells.Item(77, 11).Value = 18524517.5  # K77: 10297094 -> 18524517.5
$ws.Cells.Item(77, 12).Value = 69451935  # L77: 69451940 -> 69451935
$ws.Cells.Item(77, 13).Value = -18520149.5  # M77: -10292726 -> -18520149.5
$ws.Cells.Item(77, 14).Value = -69460671  # N77: -69460676 -> -69460671

# Row 106
$ws.Cells.Item(106, 8).Value = 0  # H106: 2370 -> 0
$ws.Cells.Item(106, 10).Value = 0  # J106: 2370 -> 0
$ws.Cells.Item(106, 12).Value = 0  # L106: 2370 -> 0
$ws.Cells.Item(106, 14).ClearContents()  # N106: -4894 -> (removed)

# Row 132
$ws.Cells.Item(132, 8).Value = 40007100  # H132: 40007224 -> 40007100
$ws.Cells.Item(132, 9).Value = 5818.3687  # I132: 5823.579 -> 5818.3687
$ws.Cells.Item(132, 10).Value = 166677840  # J132: 166678340 -> 166677840
$ws.Cells.Item(132, 11).Value = 17455.1061  # K132: 17470.737 -> 17455.1061
$ws.Cells.Item(132, 12).Value = 500033520  # L132: 500035020 -> 500033520
$ws.Cells.Item(132, 13).Value = -14925.1061  # M132: -14940.737 -> -14925.1061
$ws.Cells.Item(132, 14).Value = -500038580  # N132: -500040080 -> -500038580

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Cells.Item(105, 8).Value = 4831  # H105: 4663 -> 4831
$ws.Cells.Item(105, 9).Value = 4661.4443  # I105: 4447.909 -> 4661.4443
$ws.Cells.Item(105, 11).Value = 4661.4443  # K105: 4447.909 -> 4661.4443
$ws.Cells.Item(105, 13).Value = -2914.4443  # M105: -2700.909 -> -2914.4443

$ws = $wb.Worksheets.Item("CRP")
# Row 18
$ws.Cells.Item(18, 8).Value = 28386.75  # H18: 28387.25 -> 28386.75
$ws.Cells.Item(18, 10).Value = 28386.75  # J18: 28387.25 -> 28386.75
$ws.Cells.Item(18, 12).Value = 28386.75  # L18: 28387.25 -> 28386.75
$ws.Cells.Item(18, 14).Value = -28846.75  # N18: -28847.25 -> -28846.75

# Row 22
$ws.Cells.Item(22, 8).Value = 1206.6842  # H22: 1398.5625 -> 1206.6842
$ws.Cells.Item(22, 9).Value = 299.77777  # I22: 335.42856 -> 299.77777
$ws.Cells.Item(22, 10).Value = 2022.9  # J22: 2225.4443 -> 2022.9
$ws.Cells.Item(22, 11).Value = 299.77777  # K22: 335.42856 -> 299.77777
$ws.Cells.Item(22, 12).Value = 2022.9  # L22: 2225.4443 -> 2022.9
$ws.Cells.Item(22, 13).Value = 50.22223000000002  # M22: 14.57144 -> 50.22223000000002
$ws.Cells.Item(22, 14).Value = -2722.9  # N22: -2925.4443 -> -2722.9

# Row 31
$ws.Cells.Item(31, 8).Value = 2166.7222  # H31: 2110.5789 -> 2166.7222
$ws.Cells.Item(31, 10).Value = 4590.5  # J31: 4091.8572 -> 4590.5
$ws.Cells.Item(31, 12).Value = 4590.5  # L31: 4091.8572 -> 4590.5
$ws.Cells.Item(31, 14).Value = -5180.5  # N31: -4681.8572 -> -5180.5

# Row 34
$ws.Cells.Item(34, 8).Value = 2166.7222  # H34: 2110.5789 -> 2166.7222
$ws.Cells.Item(34, 10).Value = 4590.5  # J34: 4091.8572 -> 4590.5
$ws.Cells.Item(34, 12).Value = 4590.5  # L34: 4091.8572 -> 4590.5
$ws.Cells.Item(34, 14).Value = -4994.5  # N34: -4495.8572 -> -4994.5

# Row 58
$ws.Cells.Item(58, 8).Value = 2796.7273  # H58: 2896.5 -> 2796.7273
$ws.Cells.Item(58, 9).Value = 2650.1428  # I58: 2792 -> 2650.1428
$ws.Cells.Item(58, 11).Value = 2650.1428  # K58: 2792 -> 2650.1428
$ws.Cells.Item(58, 13).Value = -2447.1428  # M58: -2589 -> -2447.1428

# Row 99
$ws.Cells.Item(99, 8).Value = 2505.3333  # H99: 2507.25 -> 2505.3333
$ws.Cells.Item(99, 10).Value = 2514  # J99: 2513.5 -> 2514
$ws.Cells.Item(99, 12).Value = 2514  # L99: 2513.5 -> 2514
$ws.Cells.Item(99, 14).Value = -5510  # N99: -5509.5 -> -5510

# Row 103
$ws.Cells.Item(103, 8).Value = 15571  # H103: 12999.25 -> 15571
$ws.Cells.Item(103, 9).Value = 13166.167  # I103: 12999.25 -> 13166.167
$ws.Cells.Item(103, 10).Value = 30000  # J103: 0 -> 30000
$ws.Cells.Item(103, 11).Value = 13166.167  # K103: 12999.25 -> 13166.167
$ws.Cells.Item(103, 12).Value = 30000  # L103: 0 -> 30000
$ws.Cells.Item(103, 13).Value = -11994.167  # M103: -11827.25 -> -11994.167
$ws.Cells.Item(103, 14).Value = -32344  # N103: None -> -32344

# Row 105
$ws.Cells.Item(105, 8).Value = 2135.476  # H105: 1833.1852 -> 2135.476
$ws.Cells.Item(105, 9).Value = 1740.2222  # I105: 1518.8636 -> 1740.2222
$ws.Cells.Item(105, 10).Value = 4507  # J105: 3216.2 -> 4507
$ws.Cells.Item(105, 11).Value = 1740.2222  # K105: 1518.8636 -> 1740.2222
$ws.Cells.Item(105, 12).Value = 4507  # L105: 3216.2 -> 4507
$ws.Cells.Item(105, 13).Value = 6.77780000000007  # M105: 228.1364000000001 -> 6.77780000000007
$ws.Cells.Item(105, 14).Value = -8001  # N105: -6710.2 -> -8001

# Row 126
$ws.Cells.Item(126, 8).Value = 2505.3333  # H126: 2507.25 -> 2505.3333
$ws.Cells.Item(126, 10).Value = 2514  # J126: 2513.5 -> 2514
$ws.Cells.Item(126, 12).Value = 7542  # L126: 7540.5 -> 7542
$ws.Cells.Item(126, 14).Value = -12482  # N126: -12480.5 -> -12482

# Row 132
$ws.Cells.Item(132, 8).Value = 4503.077  # H132: 4388.357 -> 4503.077
$ws.Cells.Item(132, 9).Value = 4503.077  # I132: 4388.357 -> 4503.077
$ws.Cells.Item(132, 11).Value = 13509.231  # K132: 13165.071 -> 13509.231
$ws.Cells.Item(132, 13).Value = -10979.231  # M132: -10635.071 -> -10979.231

# Row 134
$ws.Cells.Item(134, 8).Value = 10003390  # H134: 11114600 -> 10003390
$ws.Cells.Item(134, 9).Value = 2987.375  # I134: 3057 -> 2987.375
$ws.Cells.Item(134, 11).Value = 8962.125  # K134: 9171 -> 8962.125
$ws.Cells.Item(134, 13).Value = -6427.125  # M134: -6636 -> -6427.125

# Row 136
$ws.Cells.Item(136, 8).Value = 2796.7273  # H136: 2896.5 -> 2796.7273
$ws.Cells.Item(136, 9).Value = 2650.1428  # I136: 2792 -> 2650.1428
$ws.Cells.Item(136, 11).Value = 7950.428400000001  # K136: 8376 -> 7950.428400000001
$ws.Cells.Item(136, 13).Value = -5400.428400000001  # M136: -5826 -> -5400.428400000001

$ws = $wb.Worksheets.Item("GSM")
# Row 129
$ws.Cells.Item(129, 8).Value = 78000  # H129: 129000 -> 78000
$ws.Cells.Item(129, 10).Value = 78000  # J129: 129000 -> 78000
$ws.Cells.Item(129, 12).Value = 78000  # L129: 129000 -> 78000
$ws.Cells.Item(129, 14).Value = -88000  # N129: -139000 -> -88000

# Row 132
$ws.Cells.Item(132, 8).Value = 1000  # H132: 2000 -> 1000
$ws.Cells.Item(132, 9).Value = 1000  # I132: 2000 -> 1000
$ws.Cells.Item(132, 11).Value = 3000  # K132: 6000 -> 3000
$ws.Cells.Item(132, 13).Value = -470  # M132: -3470 -> -470

# Row 136
$ws.Cells.Item(136, 8).Value = 41028.25  # H136: 38801 -> 41028.25
$ws.Cells.Item(136, 10).Value = 41028.25  # J136: 38801 -> 41028.25
$ws.Cells.Item(136, 12).Value = 123084.75  # L136: 116403 -> 123084.75
$ws.Cells.Item(136, 14).Value = -128184.75  # N136: -121503 -> -128184.75

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 4843.385  # H7: 5084.0586 -> 4843.385
$ws.Cells.Item(7, 9).Value = 4826.222  # I7: 5174.6665 -> 4826.222
$ws.Cells.Item(7, 10).Value = 4882  # J7: 4866.6 -> 4882
$ws.Cells.Item(7, 11).Value = 4826.222  # K7: 5174.6665 -> 4826.222
$ws.Cells.Item(7, 12).Value = 4882  # L7: 4866.6 -> 4882
$ws.Cells.Item(7, 13).Value = -4714.222  # M7: -5062.6665 -> -4714.222
$ws.Cells.Item(7, 14).Value = -5106  # N7: -5090.6 -> -5106

# Row 9
$ws.Cells.Item(9, 8).Value = 2842.6  # H9: 2805.0833 -> 2842.6
$ws.Cells.Item(9, 9).Value = 589.6667  # I9: 543.6 -> 589.6667
$ws.Cells.Item(9, 10).Value = 3808.1428  # J9: 4420.4287 -> 3808.1428
$ws.Cells.Item(9, 11).Value = 589.6667  # K9: 543.6 -> 589.6667
$ws.Cells.Item(9, 12).Value = 3808.1428  # L9: 4420.4287 -> 3808.1428
$ws.Cells.Item(9, 13).Value = -365.6667  # M9: -319.6 -> -365.6667
$ws.Cells.Item(9, 14).Value = -4256.1428  # N9: -4868.4287 -> -4256.1428

# Row 16
$ws.Cells.Item(16, 8).Value = 2288.2  # H16: 1836.7142 -> 2288.2
$ws.Cells.Item(16, 9).Value = 2288.2  # I16: 1836.7142 -> 2288.2
$ws.Cells.Item(16, 11).Value = 2288.2  # K16: 1836.7142 -> 2288.2
$ws.Cells.Item(16, 13).Value = -2118.2  # M16: -1666.7142 -> -2118.2

# Row 21
$ws.Cells.Item(21, 8).Value = 18598.6  # H21: 17361 -> 18598.6
$ws.Cells.Item(21, 9).Value = 0  # I21: 4985 -> 0
$ws.Cells.Item(21, 11).Value = 0  # K21: 4985 -> 0
$ws.Cells.Item(21, 13).ClearContents()  # M21: -4811 -> (removed)

# Row 25
$ws.Cells.Item(25, 8).Value = 16076.308  # H25: 12886.75 -> 16076.308
$ws.Cells.Item(25, 9).Value = 15000  # I25: 0 -> 15000
$ws.Cells.Item(25, 10).Value = 17332  # J25: 12886.75 -> 17332
$ws.Cells.Item(25, 11).Value = 15000  # K25: 0 -> 15000
$ws.Cells.Item(25, 12).Value = 17332  # L25: 12886.75 -> 17332
$ws.Cells.Item(25, 13).Value = -14770  # M25: None -> -14770
$ws.Cells.Item(25, 14).Value = -17792  # N25: -13346.75 -> -17792

# Row 55
$ws.Cells.Item(55, 8).Value = 1740.4166  # H55: 1734.7142 -> 1740.4166
$ws.Cells.Item(55, 9).Value = 812.125  # I55: 989.8 -> 812.125
$ws.Cells.Item(55, 11).Value = 812.125  # K55: 989.8 -> 812.125
$ws.Cells.Item(55, 13).Value = -639.125  # M55: -816.8 -> -639.125

# Row 68
$ws.Cells.Item(68, 8).Value = 2468.125  # H68: 2757 -> 2468.125
$ws.Cells.Item(68, 9).Value = 2420.4285  # I68: 2748 -> 2420.4285
$ws.Cells.Item(68, 11).Value = 2420.4285  # K68: 2748 -> 2420.4285
$ws.Cells.Item(68, 13).Value = -1671.4285  # M68: -1999 -> -1671.4285

# Row 71
$ws.Cells.Item(71, 8).Value = 2468.125  # H71: 2757 -> 2468.125
$ws.Cells.Item(71, 9).Value = 2420.4285  # I71: 2748 -> 2420.4285
$ws.Cells.Item(71, 11).Value = 12102.1425  # K71: 13740 -> 12102.1425
$ws.Cells.Item(71, 13).Value = -8358.1425  # M71: -9996 -> -8358.1425

# Row 126
$ws.Cells.Item(126, 8).Value = 4843.385  # H126: 5084.0586 -> 4843.385
$ws.Cells.Item(126, 9).Value = 4826.222  # I126: 5174.6665 -> 4826.222
$ws.Cells.Item(126, 10).Value = 4882  # J126: 4866.6 -> 4882
$ws.Cells.Item(126, 11).Value = 14478.666  # K126: 15523.9995 -> 14478.666
$ws.Cells.Item(126, 12).Value = 14646  # L126: 14599.8 -> 14646
$ws.Cells.Item(126, 13).Value = -12008.666  # M126: -13053.9995 -> -12008.666
$ws.Cells.Item(126, 14).Value = -19586  # N126: -19539.8 -> -19586

# Row 136
$ws.Cells.Item(136, 8).Value = 66670750  # H136: 66670800 -> 66670750
$ws.Cells.Item(136, 10).Value = 142862020  # J136: 142862110 -> 142862020
$ws.Cells.Item(136, 12).Value = 428586060  # L136: 428586330 -> 428586060
$ws.Cells.Item(136, 14).Value = -428591160  # N136: -428591430 -> -428591160

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Cells.Item(45, 8).Value = 17833.875  # H45: 16124.429 -> 17833.875
$ws.Cells.Item(45, 10).Value = 25596.334  # J45: 23494.5 -> 25596.334
$ws.Cells.Item(45, 12).Value = 25596.334  # L45: 23494.5 -> 25596.334
$ws.Cells.Item(45, 14).Value = -26578.334  # N45: -24476.5 -> -26578.334

# Row 81
$ws.Cells.Item(81, 8).Value = 12129.917  # H81: 17505.363 -> 12129.917
$ws.Cells.Item(81, 9).Value = 10505.363  # I81: 11455.9 -> 10505.363
$ws.Cells.Item(81, 10).Value = 30000  # J81: 78000 -> 30000
$ws.Cells.Item(81, 11).Value = 21010.726  # K81: 22911.8 -> 21010.726
$ws.Cells.Item(81, 12).Value = 60000  # L81: 156000 -> 60000
$ws.Cells.Item(81, 13).Value = -19949.726  # M81: -21850.8 -> -19949.726
$ws.Cells.Item(81, 14).Value = -62122  # N81: -158122 -> -62122

# Row 84
$ws.Cells.Item(84, 8).Value = 12129.917  # H84: 17505.363 -> 12129.917
$ws.Cells.Item(84, 9).Value = 10505.363  # I84: 11455.9 -> 10505.363
$ws.Cells.Item(84, 10).Value = 30000  # J84: 78000 -> 30000
$ws.Cells.Item(84, 11).Value = 105053.63  # K84: 114559 -> 105053.63
$ws.Cells.Item(84, 12).Value = 300000  # L84: 780000 -> 300000
$ws.Cells.Item(84, 13).Value = -99749.62999999999  # M84: -109255 -> -99749.62999999999
$ws.Cells.Item(84, 14).Value = -310608  # N84: -790608 -> -310608

# Row 92
$ws.Cells.Item(92, 8).Value = 76024.25  # H92: 81218.2 -> 76024.25
$ws.Cells.Item(92, 10).Value = 76024.25  # J92: 81218.2 -> 76024.25
$ws.Cells.Item(92, 12).Value = 76024.25  # L92: 81218.2 -> 76024.25
$ws.Cells.Item(92, 14).Value = -81016.25  # N92: -86210.2 -> -81016.25

# Row 115
$ws.Cells.Item(115, 8).Value = 24681.75  # H115: 24587.75 -> 24681.75
$ws.Cells.Item(115, 10).Value = 24681.75  # J115: 24587.75 -> 24681.75
$ws.Cells.Item(115, 12).Value = 24681.75  # L115: 24587.75 -> 24681.75
$ws.Cells.Item(115, 14).Value = -27815.75  # N115: -27721.75 -> -27815.75

# Row 121
$ws.Cells.Item(121, 8).Value = 49708.5  # H121: 59998.5 -> 49708.5
$ws.Cells.Item(121, 10).Value = 49708.5  # J121: 59998.5 -> 49708.5
$ws.Cells.Item(121, 12).Value = 49708.5  # L121: 59998.5 -> 49708.5
$ws.Cells.Item(121, 14).Value = -53202.5  # N121: -63492.5 -> -53202.5

# Row 132
$ws.Cells.Item(132, 8).Value = 1241.3529  # H132: 1484.6923 -> 1241.3529
$ws.Cells.Item(132, 9).Value = 1231.4375  # I132: 1554.6364 -> 1231.4375
$ws.Cells.Item(132, 10).Value = 1400  # J132: 1100 -> 1400
$ws.Cells.Item(132, 11).Value = 3694.3125  # K132: 4663.9092 -> 3694.3125
$ws.Cells.Item(132, 12).Value = 4200  # L132: 3300 -> 4200
$ws.Cells.Item(132, 13).Value = -1164.3125  # M132: -2133.9092 -> -1164.3125
$ws.Cells.Item(132, 14).Value = -9260  # N132: -8360 -> -9260

# Row 136
$ws.Cells.Item(136, 8).Value = 1903.6428  # H136: 1682.5883 -> 1903.6428
$ws.Cells.Item(136, 9).Value = 1473  # I136: 1267.5 -> 1473
$ws.Cells.Item(136, 11).Value = 4419  # K136: 3802.5 -> 4419
$ws.Cells.Item(136, 13).Value = -1869  # M136: -1252.5 -> -1869
